# Apply the updates described by the diff to the "Metadata" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$ws.Cells.Item(3, 2).Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$ws.Cells.Item(8, 2).Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank, now "Alvearie Team"
$ws.Cells.Item(9, 2).Value = "Alvearie Team"

# Row 10 ("Contact" / "No display for ContactDetail") becomes the new
# "Jurisdiction" / "United States of America" row.
$ws.Cells.Item(10, 1).Value = "Jurisdiction"
$ws.Cells.Item(10, 2).Value = "United States of America"

# Row 11 duplicated row 10's old "Contact" / "No display for ContactDetail"
# content and is removed entirely, shifting everything below up by one row
# (new dimension becomes A1:B19).
$ws.Rows.Item(11).Delete()
